{"js": "// Word turned \"Don't check spelling or grammar\" ON for the runs that\n// hold the five inline screenshots (this is what Word silently stamps\n// as <w:noProof/> on a picture run when it re-renders the layout), and\n// the author typed a trailing \" test\" note in a brand-new paragraph at\n// the very end of the document before submitting.\n\n// 1) Mark every inline picture's run as \"noProof\" (w:rPr/w:noProof).\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pictures.items.length; i++) {\n  const pictureRange = pictures.items[i].getRange();\n  pictureRange.hasNoProofing = true;\n}\nawait context.sync();\n\n// 2) Append a brand-new, plain paragraph (no inherited list/style)\n//    containing \" test\" right before the end of the document body.\nconst body = context.document.body;\nconst trailingParagraph = body.insertParagraph(\"\", Word.InsertLocation.end);\nconst trailingOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\"> test</w:t></w:r></w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\nconst trailingRange = trailingParagraph.getRange(Word.RangeLocation.whole);\ntrailingRange.insertOoxml(trailingOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word turned \"Don't check spelling or grammar\" ON for the runs that\n# hold the five inline screenshots (Word stamps this as <w:noProof/>\n# on a picture run when it re-renders the layout), and the author\n# typed a trailing \" test\" note in a brand-new paragraph at the very\n# end of the document right before submitting.\n\n$d = $word.ActiveDocument\n\n# 1) Mark every inline picture's run as \"noProof\".\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $shape = $d.InlineShapes.Item($i)\n    $shape.Range.NoProofing = $true\n}\n\n# 2) Append a brand-new, plain paragraph (no inherited list/style)\n#    containing \" test\" right after the last paragraph in the body.\n$lastPara = $d.Paragraphs.Last\n$endOfDoc = $lastPara.Range.End\n$insertionPoint = $d.Range($endOfDoc, $endOfDoc)\n$trailingXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\"> test</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertionPoint.InsertXML($trailingXml)\n"}
